# Insert two new rows at the top of the daily-price subset (row 520),
# pushing all existing rows from 520-622 down to 522-624, and populate
# the two newly inserted rows with the new "Bola 8" / "Sin especificar"
# observations dated 2023-11-28 (serial 45258).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("520:521").Insert()

# New row 520
$ws.Cells.Item(520, 1).Value  = 8
$ws.Cells.Item(520, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(520, 3).Value  = "Coquimbo"
$ws.Cells.Item(520, 4).Value  = 45258
$ws.Cells.Item(520, 5).Value  = 4
$ws.Cells.Item(520, 6).Value  = 100112032
$ws.Cells.Item(520, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(520, 8).Value  = "Bola 8"
$ws.Cells.Item(520, 9).Value  = "Primera"
$ws.Cells.Item(520, 10).Value = 400
$ws.Cells.Item(520, 11).Value = 11000
$ws.Cells.Item(520, 12).Value = 12000
$ws.Cells.Item(520, 13).Value = 11500
$ws.Cells.Item(520, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(520, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(520, 16).Value = 230
$ws.Cells.Item(520, 17).Value = 50
$ws.Cells.Item(520, 18).Value = "Hortaliza"

# New row 521
$ws.Cells.Item(521, 1).Value  = 8
$ws.Cells.Item(521, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(521, 3).Value  = "Coquimbo"
$ws.Cells.Item(521, 4).Value  = 45258
$ws.Cells.Item(521, 5).Value  = 4
$ws.Cells.Item(521, 6).Value  = 100112032
$ws.Cells.Item(521, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(521, 8).Value  = "Sin especificar"
$ws.Cells.Item(521, 9).Value  = "Primera"
$ws.Cells.Item(521, 10).Value = 500
$ws.Cells.Item(521, 11).Value = 10000
$ws.Cells.Item(521, 12).Value = 11000
$ws.Cells.Item(521, 13).Value = 10500
$ws.Cells.Item(521, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(521, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(521, 16).Value = 175
$ws.Cells.Item(521, 17).Value = 60
$ws.Cells.Item(521, 18).Value = "Hortaliza"
